# Quarterly financials update for CLR:
#  - Insert two new columns (new period columns) before column D, shifting
#    the existing D:K quarter columns to F:M.
#  - Populate the two new columns (D, E) with the newest two quarters of data.
#  - Fix up a handful of "Capital Expenditures" (row 91) values that were
#    revised, not just shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank columns at D:E, pushing old D:K -> F:M ---------
$ws.Columns("D:E").Insert()

# The inserted columns pick up a default/blank style. Copy number formats
# from column F (the first of the old, shifted-right columns) into the two
# new columns so each row's D/E cells match the rest of that row visually
# (date format for header rows, number format for data rows, etc.)
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill the new D & E columns with the new quarter's data ---------
$newData = @(
    @{Row=7; D=43465; E=43373},
    @{Row=8; D=1149300; E=1282200},
    @{Row=9; D=198900; E=206800},
    @{Row=10; D=950400; E=1075400},
    @{Row=12; D=3300; E=2300},
    @{Row=13; D=0; E=0},
    @{Row=14; D=30100; E=29400},
    @{Row=15; D=488400; E=469300},
    @{Row=17; D=818900; E=798000},
    @{Row=18; D=330400; E=484200},
    @{Row=20; D=1000; E=800},
    @{Row=21; D=819900; E=953500},
    @{Row=22; D=69400; E=73400},
    @{Row=23; D=262000; E=411600},
    @{Row=24; D=62900; E=97500},
    @{Row=25; D=0; E=0},
    @{Row=26; D=199100; E=314200},
    @{Row=27; D=197700; E=314200},
    @{Row=28; D=0; E=0},
    @{Row=29; D=0; E="NA"},
    @{Row=30; D=0; E=0},
    @{Row=31; D=0; E=0},
    @{Row=32; D=-1000; E=-800},
    @{Row=33; D=197700; E=314200},
    @{Row=34; D=0; E=0},
    @{Row=35; D=197700; E=314200},
    @{Row=38; D=43465; E=43373},
    @{Row=41; D=282700; E=12900},
    @{Row=42; D=0; E=0},
    @{Row=43; D=1012400; E=1237400},
    @{Row=44; D=88500; E=104200},
    @{Row=45; D=28700; E=14700},
    @{Row=46; D=1412400; E=1369100},
    @{Row=47; D=0; E=0},
    @{Row=48; D=13869800; E=13644500},
    @{Row=49; D=0; E=0},
    @{Row=50; D=0; E=0},
    @{Row=51; D=0; E=0},
    @{Row=52; D=15800; E=17400},
    @{Row=53; D=0; E=0},
    @{Row=54; D=15297900; E=15031100},
    @{Row=57; D=717600; E=771400},
    @{Row=58; D=2400; E=2300},
    @{Row=59; D=667600; E=716700},
    @{Row=60; D=1387500; E=1490400},
    @{Row=61; D=5766000; E=5955300},
    @{Row=62; D=1722600; E=1646500},
    @{Row=63; D=0; E=0},
    @{Row=64; D=0; E=0},
    @{Row=65; D=0; E=0},
    @{Row=66; D=9152800; E=9092300},
    @{Row=68; D=0; E=0},
    @{Row=69; D=0; E=0},
    @{Row=70; D=0; E=0},
    @{Row=71; D=0; E=0},
    @{Row=72; D=4706100; E=4508400},
    @{Row=73; D=0; E=0},
    @{Row=74; D=0; E=0},
    @{Row=75; D=0; E=0},
    @{Row=76; D=6145100; E=5938800},
    @{Row=77; D=0; E=0},
    @{Row=80; D=43465; E=43373},
    @{Row=81; D=197700; E=314200},
    @{Row=83; D=488400; E=468500},
    @{Row=84; D=0; E=0},
    @{Row=85; D=0; E=0},
    @{Row=86; D=0; E=0},
    @{Row=87; D=0; E=0},
    @{Row=88; D=0; E=0},
    @{Row=89; D=955300; E=860700},
    @{Row=91; D=-15900; E=-3500},
    @{Row=92; D=0; E=0},
    @{Row=93; D=0; E=0},
    @{Row=94; D=-756700; E=-759900},
    @{Row=96; D=0; E=0},
    @{Row=97; D=0; E=0},
    @{Row=98; D=0; E=0},
    @{Row=99; D=0; E=0},
    @{Row=100; D=71300; E=-218000},
    @{Row=101; D=0; E=0},
    @{Row=102; D=269900; E=-117100}
)

foreach ($entry in $newData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
}

# --- 3. Row 91 ("Capital Expenditures") got real data revisions, not just
#        a shift -- columns F:J differ from the plain old-D:H shift. -----
$ws.Cells.Item(91, 6).Value  = -9600
$ws.Cells.Item(91, 7).Value  = -2600
$ws.Cells.Item(91, 8).Value  = -5000
$ws.Cells.Item(91, 9).Value  = -2700
$ws.Cells.Item(91, 10).Value = -700

Write-Output "done"
